$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C..N, rows 2..13 (row -> column letter -> value)
$data = @{
    2  = @{ C = 54.114;               D = 54.114;               E = 3.16948242;            F = 0.00135841;            G = 0.07295569;            H = 3.97650341;            I = 6.551482727285933;   J = 6.551482727285933;   K = 0.3907925152580654;  L = 0.0001665157632159772; M = 0.008893789057972293;  N = 0.8448462905931474 }
    3  = @{ C = 90.295;               D = 90.295;               E = 1.904231;              F = 0.000943;              G = 0.08433705000000001;   H = 7.67577954;            I = 11.89283037861107;   J = 11.89283037861107;   K = 0.2527369708595069;  L = 0.0001822350394671936; M = 0.01559989991233196;   N = 2.005840126932154 }
    4  = @{ C = 27.812;               D = 55.605;               E = 3.104128;              F = 0.00264691;            G = 0.0365285;             H = 1.03036605;            I = 4.142924066383527;   J = 8.279222217124053;   K = 0.4427500415151439;  L = 0.0005123196495580892; M = 0.007658787474240003;  N = 0.3245784489303126 }
    5  = @{ C = 46.208;               D = 90.684;               E = 1.91382706;            F = 0.00166198;            G = 0.03787804;            H = 1.78094815;            I = 8.199632306998282;   J = 15.03117028311768;   K = 0.3118262287250341;  L = 0.0002498441548279681; M = 0.006012809795764376;  N = 0.5670590400082981 }
    6  = @{ C = 13.787;               D = 55.081;               E = 3.21703234;            F = 0.00435844;            G = 0.01488329;            H = 0.21439266;            I = 3.259794051356198;   J = 13.04640362887979;   K = 0.656420880163054;   L = 0.0006457398441681808; M = 0.003572530258563019;  N = 0.1058296747381356 }
    7  = @{ C = 24.351;               D = 90.449;               E = 1.92385209;            F = 0.00313757;            G = 0.01886938;            H = 0.4741455;             I = 5.201043150162268;   J = 15.66961824474151;   K = 0.3288887290705498;  L = 0.000655178642931952;  M = 0.004816035507926994;  N = 0.2096587577418312 }
    8  = @{ C = 9.026;                D = 54.047;               E = 3.2465932;             F = 0.00649477;            G = 0.00973795;            H = 0.09139376;            I = 1.854929854405488;   J = 11.1118073805868;    K = 0.6145973977206909;  L = 0.00153800201472472;   M = 0.002884756439843269;  N = 0.04736172435623935 }
    9  = @{ C = 16.821;               D = 86.29600000000001;    E = 2.01509231;            F = 0.004107620000000001;  G = 0.01131171;            H = 0.19793911;            I = 4.211020267309538;   J = 14.80106264615671;   K = 0.3395036868161016;  L = 0.000663986455896213;  M = 0.00246834225531821;   N = 0.09679939926366127 }
    10 = @{ C = 6.91;                 D = 55.089;               E = 3.19726144;            F = 0.00728126;            G = 0.00632601;            H = 0.04580960000000001;   I = 1.412053077576723;   J = 11.23250854654686;   K = 0.662175724919601;   L = 0.0009135220495067073; M = 0.001703840324605493;  N = 0.02170012876063707 }
    11 = @{ C = 12.66;                D = 80.422;               E = 2.16211396;            F = 0.005118800000000001;  G = 0.00795972;            H = 0.10568401;            I = 3.356138206877611;   J = 13.83876538667476;   K = 0.3621626749625843;  L = 0.0008798775736115477; M = 0.001979683931861364;  N = 0.05539014756133772 }
    12 = @{ C = 5.608;                D = 55.843;               E = 3.180136210000001;     F = 0.00949055;            G = 0.005367089999999999;  H = 0.03192714999999999;   I = 1.275916911073744;   J = 12.76625741273958;   K = 0.7060770387562498;  L = 0.001716747628488218;  M = 0.001721453174772097;  N = 0.01763455825424392 }
    13 = @{ C = 10.617;               D = 75.434;               E = 2.31641947;            F = 0.00557833;            G = 0.00581731;            H = 0.06624899999999999;   I = 3.315615222595092;   J = 14.03512125450395;   K = 0.4230620195136572;  L = 0.0009005373276986239; M = 0.001655029542541971;  N = 0.04440052480999975 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
